$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Student IDs to populate into the new column I (rows 2-12), mirroring the
# formatting already used in column B of the same row.
$studentIds = @{
    2  = "17-0025"
    3  = "18-0186"
    4  = "18-0146"
    5  = "17-0052"
    6  = "17-0029"
    7  = "18-0074"
    8  = "17-0076"
    9  = "17-0072"
    10 = "18-0095"
    11 = "17-0090"
    12 = "17-0005"
}

foreach ($row in 2..12) {
    $bCell = $ws.Range("B" + $row)
    $iCell = $ws.Range("I" + $row)
    $bCell.Copy()
    $iCell.PasteSpecial(-4122)
    $iCell.Value = $studentIds[$row]
}

# Conditional formatting on I9 ("student_id" column) highlighting cells that
# contain "16-" or "17-" prefixes (three rules, matching the source commit).
$fc = $ws.Range("I9").FormatConditions

$cond1 = $fc.Add(9, 0, "16-")
$cond1.Text = "16-"
$cond1.Priority = 2
$cond1.Font.Color = 0

$cond2 = $fc.Add(9, 0, "17-")
$cond2.Text = "17-"
$cond2.Priority = 3
$cond2.Font.Color = 0

$cond3 = $fc.Add(9, 0, "16-")
$cond3.Text = "16-"
$cond3.Priority = 4
$cond3.Font.Color = 0

# Scroll/zoom/selection state left behind after the edit.
$excel.ActiveWindow.DisplayGridlines = $true
$excel.ActiveWindow.Zoom = 90
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("I2:I12").Select()
